$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.129.51"
$ws.Range("E2").Value = "  -4.10%  "
$ws.Range("D3").Value = "3.240.77"
$ws.Range("E3").Value = "  -5.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "175.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "519.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.593"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.78%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.246.32"
$ws.Range("E8").Value = "  -4.86%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.602"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.131"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.58%  "
$ws.Range("D15").Value = "3.743.99"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.115"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.62%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.229.99"
$ws.Range("E17").Value = "  -5.48%  "
$ws.Range("D18").Value = "62.957.02"
$ws.Range("E18").Value = "  -3.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.960"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "367.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "654.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.86%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.374"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0713"
$ws.Range("E41").Value = "  +12.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.123"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.51%  "
$ws.Range("D43").Value = "2.865.42"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0390"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("E47").Value = "  -7.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  -2.94%  "
